# Flip the "BIEfIE" boolean control-lever value (cell B2 on sheet "BIEfIE")
# from 1 (include emissions from imported electricity) to 0 (exclude them).
$wb = $excel.ActiveWorkbook

$leverSheet = $wb.Worksheets.Item("BIEfIE")
$leverSheet.Range("B2").Value = 0

# Leave the "About" sheet as the selected/active sheet in the saved view.
$aboutSheet = $wb.Worksheets.Item("About")
$aboutSheet.Activate()
